$d = $word.ActiveDocument

# --- 1) Merge the "1." + "5" + ".1" runs in the version-number cell into a
#        single run reading "1.5.1" (the three runs already share the same
#        run formatting, so Find/Replace naturally coalesces them). ---
$d.Content.Find.Execute("1.5.1", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1.5.1", 2) | Out-Null

# --- 2) Merge the "3 สิงหาคม" + " 2564" runs into a single run reading
#        "3 สิงหาคม 2564". ---
$d.Content.Find.Execute("3 สิงหาคม 2564", $true, $false, $false, $false, $false, `
    $true, 1, $false, "3 สิงหาคม 2564", 2) | Out-Null

# --- 3) Replace the "ผู้จัดทำ" name/role in the last table row:
#        "วริศรา (D)" -> "ปรีชญา (PM)", and while doing so drop the
#        w:hint="cs" rendering hint that the original "วริศรา " run carried
#        (the replacement text keeps w:cs but no longer hints cs). We
#        rebuild the paragraph's OOXML directly so the run formatting comes
#        out exactly right. ---
$tbl = $d.Tables.Item(1)
$nameCell = $tbl.Cell($tbl.Rows.Count, 4)
$nameCell.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="706E9FD1" w14:textId="7007F94F" w:rsidR="00761B6A" w:rsidRDefault="00761B6A" w:rsidP="00FB6638"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t xml:space="preserve">ปรีชญา </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>(PM)</w:t></w:r></w:p>
"@) | Out-Null

# --- 4) The neighbouring "กิตติพศ (SP)" cell keeps its text, but loses the
#        same w:hint="cs" rendering hint on its first run. ---
$codeCell = $tbl.Cell($tbl.Rows.Count, 5)
$codeCell.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5EF61068" w14:textId="2F7FA5CE" w:rsidR="00667B81" w:rsidRDefault="00761B6A" w:rsidP="00FB6638"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:cs/></w:rPr><w:t xml:space="preserve">กิตติพศ </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TH Sarabun New" w:hAnsi="TH Sarabun New"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>(SP)</w:t></w:r></w:p>
"@) | Out-Null
